$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 4039.5908
$ws.Range("I32").Value = 3525.6
$ws.Range("J32").Value = 4190.7646
$ws.Range("K32").Value = 3525.6
$ws.Range("L32").Value = 4190.7646
$ws.Range("M32").Value = -3199.6
$ws.Range("N32").Value = -4842.7646
$ws.Range("H64").Value = 8639.538
$ws.Range("I64").Value = 7249.75
$ws.Range("J64").Value = 10863.2
$ws.Range("K64").Value = 7249.75
$ws.Range("L64").Value = 10863.2
$ws.Range("M64").Value = -7001.75
$ws.Range("N64").Value = -11359.2
$ws.Range("H67").Value = 8639.538
$ws.Range("I67").Value = 7249.75
$ws.Range("J67").Value = 10863.2
$ws.Range("K67").Value = 7249.75
$ws.Range("L67").Value = 10863.2
$ws.Range("M67").Value = -6391.75
$ws.Range("N67").Value = -12579.2
$ws.Range("H80").Value = 500.05554
$ws.Range("I80").Value = 434.1111
$ws.Range("J80").Value = 566
$ws.Range("K80").Value = 1302.3333
$ws.Range("L80").Value = 1698
$ws.Range("M80").Value = -304.3333
$ws.Range("N80").Value = -3694
$ws.Range("H83").Value = 500.05554
$ws.Range("I83").Value = 434.1111
$ws.Range("J83").Value = 566
$ws.Range("K83").Value = 3906.9999
$ws.Range("L83").Value = 5094
$ws.Range("M83").Value = 1085.0001
$ws.Range("N83").Value = -15078
$ws.Range("H106").Value = 23414.934
$ws.Range("I106").Value = 30433.4
$ws.Range("K106").Value = 30433.4
$ws.Range("M106").Value = -29802.4
$ws.Range("H132").Value = 3326.8
$ws.Range("I132").Value = 3099.7
$ws.Range("J132").Value = 3781
$ws.Range("K132").Value = 9299.099999999999
$ws.Range("L132").Value = 11343
$ws.Range("M132").Value = -6769.099999999999
$ws.Range("N132").Value = -16403
$ws.Range("H138").Value = 15154020
$ws.Range("I138").Value = 41667996
$ws.Range("J138").Value = 3175.9285
$ws.Range("K138").Value = 125003988
$ws.Range("L138").Value = 9527.7855
$ws.Range("M138").Value = -124998848
$ws.Range("N138").Value = -19807.7855
$ws.Range("H141").Value = 3174.3635
$ws.Range("I141").Value = 2534.5789
$ws.Range("K141").Value = 7603.736699999999
$ws.Range("M141").Value = -2423.736699999999
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1767.9844
$ws.Range("I32").Value = 1680.6936
$ws.Range("K32").Value = 1680.6936
$ws.Range("M32").Value = -1393.6936
$ws.Range("H61").Value = 2897.6191
$ws.Range("I61").Value = 1499
$ws.Range("K61").Value = 1499
$ws.Range("M61").Value = -1287
$ws.Range("H74").Value = 3325.8823
$ws.Range("I74").Value = 2837.7715
$ws.Range("J74").Value = 4393.625
$ws.Range("K74").Value = 2837.7715
$ws.Range("L74").Value = 4393.625
$ws.Range("M74").Value = -1963.7715
$ws.Range("N74").Value = -6141.625
$ws.Range("H77").Value = 3325.8823
$ws.Range("I77").Value = 2837.7715
$ws.Range("J77").Value = 4393.625
$ws.Range("K77").Value = 14188.8575
$ws.Range("L77").Value = 21968.125
$ws.Range("M77").Value = -9820.857499999998
$ws.Range("N77").Value = -30704.125
$ws.Range("H110").Value = 2200
$ws.Range("I110").Value = 2200
$ws.Range("K110").Value = 2200
$ws.Range("M110").Value = -155
$ws.Range("H136").Value = 2897.6191
$ws.Range("I136").Value = 1499
$ws.Range("K136").Value = 4497
$ws.Range("M136").Value = -1947
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H103").Value = 20657
$ws.Range("J103").Value = 20657
$ws.Range("L103").Value = 20657
$ws.Range("N103").Value = -23001
$ws.Range("H105").Value = 1900
$ws.Range("I105").Value = 1900
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 1900
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = -153
$ws.Range("N105").ClearContents()
$ws.Range("H134").Value = 3853.1333
$ws.Range("I134").Value = 2654.9583
$ws.Range("K134").Value = 7964.874899999999
$ws.Range("M134").Value = -5429.874899999999
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1407.6
$ws.Range("I22").Value = 439.16666
$ws.Range("K22").Value = 439.16666
$ws.Range("M22").Value = -89.16665999999998
$ws.Range("H70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("N70").ClearContents()
$ws.Range("H73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("N73").ClearContents()
$ws.Range("H132").Value = 1490.4286
$ws.Range("I132").Value = 1421.5625
$ws.Range("K132").Value = 4264.6875
$ws.Range("M132").Value = -1734.6875
$ws.Range("H134").Value = 1931.65
$ws.Range("I134").Value = 1899.4445
$ws.Range("K134").Value = 5698.333500000001
$ws.Range("M134").Value = -3163.333500000001
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 197.33333
$ws.Range("J23").Value = 259.33334
$ws.Range("L23").Value = 778.0000200000001
$ws.Range("N23").Value = -1248.00002
$ws.Range("H24").Value = 537.25
$ws.Range("I24").Value = 650
$ws.Range("J24").Value = 424.5
$ws.Range("K24").Value = 1950
$ws.Range("L24").Value = 1273.5
$ws.Range("M24").Value = -1720
$ws.Range("N24").Value = -1733.5
$ws.Range("H33").Value = 474.41177
$ws.Range("J33").Value = 498.41666
$ws.Range("L33").Value = 2990.49996
$ws.Range("N33").Value = -3556.49996
$ws.Range("H87").Value = 39499.5
$ws.Range("I87").Value = 30000
$ws.Range("K87").Value = 90000
$ws.Range("M87").Value = -88752
$ws.Range("H90").Value = 39499.5
$ws.Range("I90").Value = 30000
$ws.Range("K90").Value = 270000
$ws.Range("M90").Value = -263760
$ws.Range("H117").Value = 51450.855
$ws.Range("J117").Value = 63454.53
$ws.Range("L117").Value = 190363.59
$ws.Range("N117").Value = -197247.59
$ws.Range("H137").Value = 1686.8334
$ws.Range("I137").Value = 835.2
$ws.Range("J137").Value = 2751.375
$ws.Range("K137").Value = 2505.6
$ws.Range("L137").Value = 8254.125
$ws.Range("M137").Value = 2594.4
$ws.Range("N137").Value = -18454.125
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3167.1482
$ws.Range("I132").Value = 3172.56
$ws.Range("J132").Value = 3099.5
$ws.Range("K132").Value = 9517.68
$ws.Range("L132").Value = 9298.5
$ws.Range("M132").Value = -6987.68
$ws.Range("N132").Value = -14358.5
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3742.1191
$ws.Range("I40").Value = 3700.1282
$ws.Range("K40").Value = 3700.1282
$ws.Range("M40").Value = -3564.1282
$ws.Range("H48").Value = 15677
$ws.Range("I48").Value = 18041
$ws.Range("J48").Value = 14495
$ws.Range("K48").Value = 18041
$ws.Range("L48").Value = 14495
$ws.Range("M48").Value = -17380
$ws.Range("N48").Value = -15817
$ws.Range("H61").Value = 1867.091
$ws.Range("I61").Value = 1505.5714
$ws.Range("K61").Value = 1505.5714
$ws.Range("M61").Value = -1303.5714
$ws.Range("H113").Value = 1867.091
$ws.Range("I113").Value = 1505.5714
$ws.Range("K113").Value = 1505.5714
$ws.Range("M113").Value = 664.4286
$ws.Range("H132").Value = 3212.41
$ws.Range("I132").Value = 3329.0833
$ws.Range("J132").Value = 2599.875
$ws.Range("K132").Value = 9987.249899999999
$ws.Range("L132").Value = 7799.625
$ws.Range("M132").Value = -7457.249899999999
$ws.Range("N132").Value = -12859.625
$ws.Range("H136").Value = 5199.6113
$ws.Range("I136").Value = 4407
$ws.Range("K136").Value = 13221
$ws.Range("M136").Value = -10671
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H29").Value = 286.66666
$ws.Range("J29").Value = 286.66666
$ws.Range("L29").Value = 286.66666
$ws.Range("N29").Value = -866.66666
$ws.Range("H33").Value = 49998
$ws.Range("J33").Value = 49998
$ws.Range("L33").Value = 49998
$ws.Range("N33").Value = -50498
$ws.Range("H36").Value = 49998
$ws.Range("J36").Value = 49998
$ws.Range("L36").Value = 49998
$ws.Range("N36").Value = -50498
$ws.Range("H40").Value = 30025
$ws.Range("I40").Value = 30025
$ws.Range("K40").Value = 30025
$ws.Range("M40").Value = -29876
$ws.Range("H47").Value = 25000
$ws.Range("J47").Value = 25000
$ws.Range("L47").Value = 25000
$ws.Range("N47").Value = -26144
$ws.Range("H81").Value = 5704.3335
$ws.Range("I81").Value = 4675.25
$ws.Range("K81").Value = 9350.5
$ws.Range("M81").Value = -8289.5
$ws.Range("H84").Value = 5704.3335
$ws.Range("I84").Value = 4675.25
$ws.Range("K84").Value = 46752.5
$ws.Range("M84").Value = -41448.5
$ws.Range("H119").Value = 76213
$ws.Range("J119").Value = 76213
$ws.Range("L119").Value = 76213
$ws.Range("N119").Value = -85889
$ws.Range("H132").Value = 881.93335
$ws.Range("I132").Value = 992.5714
$ws.Range("J132").Value = 785.125
$ws.Range("K132").Value = 2977.7142
$ws.Range("L132").Value = 2355.375
$ws.Range("M132").Value = -447.7142000000003
$ws.Range("N132").Value = -7415.375
